$d = $word.ActiveDocument

$d.Content.Find.Execute("olá", $false, $false, $false, $false, $false, $true, 1, $false, "Lorem ipsum dolor sit amet consectetur adipisicing elit.", 2) | Out-Null
$d.Content.Find.Execute("tudo bem?", $false, $false, $false, $false, $false, $true, 1, $false, "Maxime mollitia, molestiae quas", 2) | Out-Null
$d.Content.Find.Execute("e se eu te dissesse que vim do futuro?", $false, $false, $false, $false, $false, $true, 1, $false, "Impedit sit sunt quaerat, odit.", 2) | Out-Null
$d.Content.Find.Execute("vim pra te alertar sobre algumas coisas ", $false, $false, $false, $false, $false, $true, 1, $false, "Quo neque error repudiandae fuga? Ipsa laudantium molestias eos ", 2) | Out-Null
$d.Content.Find.Execute("que vão acontecer", $false, $false, $false, $false, $false, $true, 1, $false, "sapiente officiis modi at sunt excepturi expedita sint? Sed quibusdam", 2) | Out-Null
